# Added third location ("Greenwoods, Pasig City") to the GPS accuracy
# dataset and refreshed the existing location labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Refresh the two existing location labels -------------------------
# (text only changes; same two merged blocks E2:E31 and E32:E61)
$ws.Range("E2").Value  = "Pililla,  Rizal"
$ws.Range("E32").Value = "Pineda, Pasig City"

# --- 2. Append the new GPS readings for the third location ---------------
$newRows = @(
    @(14.562879000000001, 121.09813800000001, 14.562817000000001, 121.09808700000001),
    @(14.562879000000001, 121.098139, 14.562817000000001, 121.09808700000001),
    @(14.56288, 121.09814, 14.562817000000001, 121.09808700000001),
    @(14.56288, 121.09814, 14.562817000000001, 121.09808700000001),
    @(14.562881000000001, 121.09814, 14.562817000000001, 121.09808700000001),
    @(14.562882, 121.098141, 14.562817000000001, 121.09808700000001),
    @(14.562882, 121.09814, 14.562817000000001, 121.09808700000001),
    @(14.562882, 121.09814, 14.562817000000001, 121.09808700000001),
    @(14.562882999999999, 121.098139, 14.562817000000001, 121.09808700000001),
    @(14.562882999999999, 121.098139, 14.562817000000001, 121.09808700000001),
    @(14.562885, 121.098139, 14.562817000000001, 121.09808700000001),
    @(14.562886000000001, 121.09814, 14.562817000000001, 121.09808700000001),
    @(14.562886000000001, 121.09814, 14.562817000000001, 121.09808700000001),
    @(14.562887, 121.09814, 14.562817000000001, 121.09808700000001),
    @(14.562887999999999, 121.09814, 14.562817000000001, 121.09808700000001),
    @(14.562889, 121.09814, 14.562817000000001, 121.09808700000001),
    @(14.562889999999999, 121.09814, 14.562817000000001, 121.09808700000001),
    @(14.562891, 121.098141, 14.562817000000001, 121.09808700000001),
    @(14.562893000000001, 121.09814299999999, 14.562817000000001, 121.09808700000001),
    @(14.562893000000001, 121.098144, 14.562817000000001, 121.09808700000001),
    @(14.562894, 121.098144, 14.562817000000001, 121.09808700000001),
    @(14.562894, 121.098145, 14.562817000000001, 121.09808700000001),
    @(14.562894999999999, 121.098146, 14.562817000000001, 121.09808700000001),
    @(14.562894999999999, 121.09814799999999, 14.562817000000001, 121.09808700000001),
    @(14.562894999999999, 121.09815, 14.562817000000001, 121.09808700000001),
    @(14.562894999999999, 121.098151, 14.562817000000001, 121.09808700000001),
    @(14.562894999999999, 121.098151, 14.562817000000001, 121.09808700000001),
    @(14.562894999999999, 121.098151, 14.562817000000001, 121.09808700000001),
    @(14.562894999999999, 121.098151, 14.562817000000001, 121.09808700000001),
    @(14.562894, 121.098151, 14.562817000000001, 121.09808700000001)
)

$firstNewRow = 62
$lastNewRow = $firstNewRow + $newRows.Length - 1

# Copy the formatting of the last existing data row (61) down across the
# whole new block so the new cells pick up the same styles (centered/
# wrapped style "2" on column E, plain style "1" on A:D) as the rest of
# the table.
$ws.Range("A61:E61").Copy()
$ws.Range("A" + $firstNewRow + ":E" + $lastNewRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$r = $firstNewRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}

# Location label for the new block (merged E62:E91), same as a fresh
# shared-string entry holding the new place name.
$ws.Range("E" + $firstNewRow).Value = "Greenwoods, Pasig City"

# --- 3. Merge the Location column for the new block -----------------------
$ws.Range("E" + $firstNewRow + ":E" + $lastNewRow).Merge()

# --- 4. Scroll / select to match the new active area -----------------------
$ws.Range("E" + $firstNewRow + ":E" + $lastNewRow).Select()
